$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("P3").Value = 2.12
$ws.Range("S3").Value = 3
$ws.Range("T3").Value = 1.76
$ws.Range("X3").Value = 20
$ws.Range("Y3").Value = 10.5
$ws.Range("Z3").Value = 12.5
$ws.Range("AH3").Value = 18
$ws.Range("AO3").Value = 11.5

# Row 4
$ws.Range("N4").Value = 5.4
$ws.Range("P4").Value = 2.48
$ws.Range("Q4").Value = 1.61
$ws.Range("R4").Value = 1.6
$ws.Range("U4").Value = 2.62
$ws.Range("AC4").Value = 9.199999999999999

# Row 5
$ws.Range("F5").Value = 1.99
$ws.Range("G5").Value = 2.7
$ws.Range("H5").Value = 1.59
$ws.Range("I5").Value = 5.1
$ws.Range("Q5").Value = 1.62

# Row 6
$ws.Range("G6").Value = 2.5

# Row 9
$ws.Range("F9").Value = 2.02
$ws.Range("H9").Value = 4.6
$ws.Range("K9").Value = 3.35
$ws.Range("P9").Value = 1.45
$ws.Range("Q9").Value = 2.6

# Row 10
$ws.Range("F10").Value = 2.48
$ws.Range("H10").Value = 2.98
$ws.Range("I10").Value = 3.85
$ws.Range("J10").Value = 2.58
$ws.Range("K10").Value = 3.15
$ws.Range("P10").Value = 1.39
